$wb = $excel.ActiveWorkbook

$wsMeta    = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from FSIII")

# --- Bump the term version: 1.0.0 -> 1.1.0 ---
$wsMeta.Range("B3").Value = "1.1.0"

# --- Bump the term date ---
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- Fix the missing applyAlignment on the shared cell formats: the ---
# --- "vertical top / wrap text" alignment was present in the styles ---
# --- but never actually turned on, so turn it on for every data cell. ---
$wsMeta.Range("A1:B14").WrapText = $true
$wsInclude.Range("A1:B7").WrapText = $true
